$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns (B, C, D, E) remain stored as text so formats like
# trailing zeros ("1.00"), thousand-dot separators ("61.902.85") and
# percentage strings with padding spaces are preserved exactly.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.902.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.992.53"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.10"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.22"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.994.40"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.486"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.70"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +9.26%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.53%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.56%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.63"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.427.98"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.914.52"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.978.65"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.54"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.89"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.44"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.651"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.14"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.15"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.53"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.52"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.86%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.98"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.48%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.20"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.48"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "54.06"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.60%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "447.44"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0802"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0386"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.930.51"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -9.41%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.97"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.45"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.47"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.82%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.98"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "114.44"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0485"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.24"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.15%  "
